$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.255.09"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.578.46"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "2.586.39"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("E11").Value = "  +2.85%  "
$ws.Range("E12").Value = "  +11.06%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "3.030.81"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "59.264.02"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.43%  "
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("D18").Value = "2.581.65"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.463"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.874"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "296.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.34%  "
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.87%  "
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.595"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").Value = "1.956.01"
$ws.Range("E51").Value = "  +0.40%  "
